# Update "Pais" sheet (COVID-19 country stats) to the 22-May-2020 19:05 snapshot.
# Country-name cells are rewritten where the sort-by-total-cases order shifted a row;
# all other cells just get refreshed totals/new-cases/active/recovered/critical/deaths.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 22 de Mayo de 2020 a las 19:05'

# Row 4
$ws.Cells.Item(4, 2).Value = 1629176
$ws.Cells.Item(4, 3).Value = 8274
$ws.Cells.Item(4, 4).Value = 384269
$ws.Cells.Item(4, 5).Value = 1148163
$ws.Cells.Item(4, 7).Value = 390
$ws.Cells.Item(4, 8).Value = 96744

# Row 6
$ws.Cells.Item(6, 2).Value = 314769
$ws.Cells.Item(6, 3).Value = 3848
$ws.Cells.Item(6, 5).Value = 168542
$ws.Cells.Item(6, 7).Value = 185
$ws.Cells.Item(6, 8).Value = 20267

# Row 12
$ws.Cells.Item(12, 2).Value = 154500
$ws.Cells.Item(12, 3).Value = 952
$ws.Cells.Item(12, 4).Value = 116111
$ws.Cells.Item(12, 5).Value = 34113
$ws.Cells.Item(12, 7).Value = 27
$ws.Cells.Item(12, 8).Value = 4276

# Row 14
$ws.Cells.Item(14, 4).Value = 51307
$ws.Cells.Item(14, 5).Value = 69059

# Row 17
$ws.Cells.Item(17, 2).Value = 81767
$ws.Cells.Item(17, 3).Value = 443
$ws.Cells.Item(17, 4).Value = 41975
$ws.Cells.Item(17, 5).Value = 33612
$ws.Cells.Item(17, 7).Value = 28
$ws.Cells.Item(17, 8).Value = 6180

# Row 25
$ws.Cells.Item(25, 2).Value = 35828
$ws.Cells.Item(25, 3).Value = 522
$ws.Cells.Item(25, 5).Value = 29215
$ws.Cells.Item(25, 7).Value = 117
$ws.Cells.Item(25, 8).Value = 3056

# Row 44
$ws.Cells.Item(44, 2).Value = 15786
$ws.Cells.Item(44, 3).Value = 783
$ws.Cells.Item(44, 4).Value = 4374
$ws.Cells.Item(44, 5).Value = 10705
$ws.Cells.Item(44, 7).Value = 11
$ws.Cells.Item(44, 8).Value = 707

# Row 53
$ws.Cells.Item(53, 2).Value = 8770
$ws.Cells.Item(53, 3).Value = 16
$ws.Cells.Item(53, 4).Value = 6019
$ws.Cells.Item(53, 5).Value = 2439
$ws.Cells.Item(53, 7).Value = 6
$ws.Cells.Item(53, 8).Value = 312

# Row 58
$ws.Cells.Item(58, 2).Value = 7332
$ws.Cells.Item(58, 3).Value = 121
$ws.Cells.Item(58, 4).Value = 4377
$ws.Cells.Item(58, 5).Value = 2758

# Row 80: Tayikistan
$ws.Cells.Item(80, 1).Value = 'Tayikistan'
$ws.Cells.Item(80, 2).Value = 2551
$ws.Cells.Item(80, 3).Value = 201
$ws.Cells.Item(80, 4).Value = 1089
$ws.Cells.Item(80, 5).Value = 1418
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 44

# Row 81: Guatemala
$ws.Cells.Item(81, 1).Value = 'Guatemala'
$ws.Cells.Item(81, 2).Value = 2512
$ws.Cells.Item(81, 3).Value = 247
$ws.Cells.Item(81, 4).Value = 222
$ws.Cells.Item(81, 5).Value = 2242
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 48

# Row 82: Bulgaria
$ws.Cells.Item(82, 1).Value = 'Bulgaria'
$ws.Cells.Item(82, 3).Value = 41
$ws.Cells.Item(82, 4).Value = 769
$ws.Cells.Item(82, 5).Value = 1478
$ws.Cells.Item(82, 7).Value = 5
$ws.Cells.Item(82, 8).Value = 125

# Row 83: Bosnia y Herzegovina
$ws.Cells.Item(83, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(83, 2).Value = 2372
$ws.Cells.Item(83, 3).Value = 22
$ws.Cells.Item(83, 4).Value = 1614
$ws.Cells.Item(83, 5).Value = 617
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 141

# Row 85: Republica de Yibuti
$ws.Cells.Item(85, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(85, 2).Value = 2270
$ws.Cells.Item(85, 3).Value = 223
$ws.Cells.Item(85, 4).Value = 1064
$ws.Cells.Item(85, 5).Value = 1196
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 10

# Row 86: Croacia
$ws.Cells.Item(86, 1).Value = 'Croacia'
$ws.Cells.Item(86, 2).Value = 2243
$ws.Cells.Item(86, 3).Value = 6
$ws.Cells.Item(86, 4).Value = 2011
$ws.Cells.Item(86, 5).Value = 133
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 99

# Row 105: Sri Lanka
$ws.Cells.Item(105, 1).Value = 'Sri Lanka'
$ws.Cells.Item(105, 2).Value = 1068
$ws.Cells.Item(105, 3).Value = 20
$ws.Cells.Item(105, 4).Value = 620
$ws.Cells.Item(105, 5).Value = 439
$ws.Cells.Item(105, 8).Value = 9

# Row 106: Hong Kong
$ws.Cells.Item(106, 1).Value = 'Hong Kong'
$ws.Cells.Item(106, 2).Value = 1066
$ws.Cells.Item(106, 3).Value = 2
$ws.Cells.Item(106, 4).Value = 1029
$ws.Cells.Item(106, 5).Value = 33
$ws.Cells.Item(106, 8).Value = 4

# Row 117
$ws.Cells.Item(117, 2).Value = 838
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 4).Value = 264
$ws.Cells.Item(117, 5).Value = 563

# Row 118
$ws.Cells.Item(118, 2).Value = 814
$ws.Cells.Item(118, 3).Value = 2
$ws.Cells.Item(118, 4).Value = 672
$ws.Cells.Item(118, 5).Value = 90

# Row 119
$ws.Cells.Item(119, 4).Value = 652
$ws.Cells.Item(119, 5).Value = 59

# Row 133: Republica de Africa Central
$ws.Cells.Item(133, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(133, 2).Value = 479
$ws.Cells.Item(133, 3).Value = 43
$ws.Cells.Item(133, 4).Value = 18
$ws.Cells.Item(133, 5).Value = 461
$ws.Cells.Item(133, 8).Value = 0

# Row 134: Congo
$ws.Cells.Item(134, 1).Value = 'Congo'
$ws.Cells.Item(134, 2).Value = 469
$ws.Cells.Item(134, 4).Value = 137
$ws.Cells.Item(134, 5).Value = 316
$ws.Cells.Item(134, 8).Value = 16

# Row 135: Reunion
$ws.Cells.Item(135, 1).Value = 'Reunion'
$ws.Cells.Item(135, 2).Value = 449
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 411
$ws.Cells.Item(135, 5).Value = 37
$ws.Cells.Item(135, 8).Value = 1

# Row 136: Madagascar
$ws.Cells.Item(136, 1).Value = 'Madagascar'
$ws.Cells.Item(136, 2).Value = 448
$ws.Cells.Item(136, 3).Value = 43
$ws.Cells.Item(136, 4).Value = 135
$ws.Cells.Item(136, 5).Value = 311
$ws.Cells.Item(136, 8).Value = 2

# Row 137: Taiwan
$ws.Cells.Item(137, 1).Value = 'Taiwan'
$ws.Cells.Item(137, 2).Value = 441
$ws.Cells.Item(137, 4).Value = 408
$ws.Cells.Item(137, 5).Value = 26
$ws.Cells.Item(137, 8).Value = 7

# Row 151
$ws.Cells.Item(151, 2).Value = 225
$ws.Cells.Item(151, 3).Value = 5
$ws.Cells.Item(151, 4).Value = 119
$ws.Cells.Item(151, 5).Value = 104

# Row 174: Comoras
$ws.Cells.Item(174, 1).Value = 'Comoras'
$ws.Cells.Item(174, 2).Value = 78
$ws.Cells.Item(174, 3).Value = 44
$ws.Cells.Item(174, 4).Value = 18
$ws.Cells.Item(174, 5).Value = 59
$ws.Cells.Item(174, 8).Value = 1

# Row 175: San Martin (Parte Holandesa)
$ws.Cells.Item(175, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(175, 2).Value = 77
$ws.Cells.Item(175, 4).Value = 59
$ws.Cells.Item(175, 5).Value = 3
$ws.Cells.Item(175, 8).Value = 15

# Row 176: Malaui
$ws.Cells.Item(176, 1).Value = 'Malaui'
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 27
$ws.Cells.Item(176, 5).Value = 42

# Row 177: Libia
$ws.Cells.Item(177, 1).Value = 'Libia'
$ws.Cells.Item(177, 2).Value = 72
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(177, 4).Value = 38
$ws.Cells.Item(177, 5).Value = 31

# Row 178: Angola
$ws.Cells.Item(178, 1).Value = 'Angola'
$ws.Cells.Item(178, 3).Value = 2
$ws.Cells.Item(178, 4).Value = 17
$ws.Cells.Item(178, 5).Value = 40
$ws.Cells.Item(178, 8).Value = 3

# Row 179: Polinesia Francesa
$ws.Cells.Item(179, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(179, 2).Value = 60
$ws.Cells.Item(179, 4).Value = 60
$ws.Cells.Item(179, 5).Value = 0
$ws.Cells.Item(179, 8).Value = 0

# Row 180: Siria
$ws.Cells.Item(180, 1).Value = 'Siria'
$ws.Cells.Item(180, 2).Value = 58
$ws.Cells.Item(180, 4).Value = 36
$ws.Cells.Item(180, 5).Value = 19
$ws.Cells.Item(180, 8).Value = 3

# Row 181: Zimbabue
$ws.Cells.Item(181, 1).Value = 'Zimbabue'
$ws.Cells.Item(181, 2).Value = 51
$ws.Cells.Item(181, 4).Value = 18
$ws.Cells.Item(181, 5).Value = 29
$ws.Cells.Item(181, 8).Value = 4

# Row 182: Macao
$ws.Cells.Item(182, 1).Value = 'Macao'
$ws.Cells.Item(182, 2).Value = 45
$ws.Cells.Item(182, 4).Value = 45
$ws.Cells.Item(182, 5).Value = 0
$ws.Cells.Item(182, 8).Value = 0

# Row 183: Burundi
$ws.Cells.Item(183, 1).Value = 'Burundi'
$ws.Cells.Item(183, 2).Value = 42
$ws.Cells.Item(183, 4).Value = 20
$ws.Cells.Item(183, 5).Value = 21
$ws.Cells.Item(183, 8).Value = 1

# Row 184: San Martin (Parte Francesa)
$ws.Cells.Item(184, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(184, 2).Value = 40
$ws.Cells.Item(184, 4).Value = 33
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 8).Value = 3

# Row 185: Puerto Rico
$ws.Cells.Item(185, 1).Value = 'Puerto Rico'
$ws.Cells.Item(185, 4).Value = 1
$ws.Cells.Item(185, 5).Value = 36
$ws.Cells.Item(185, 8).Value = 2

# Row 186: Eritrea
$ws.Cells.Item(186, 1).Value = 'Eritrea'
$ws.Cells.Item(186, 2).Value = 39
$ws.Cells.Item(186, 4).Value = 39
$ws.Cells.Item(186, 5).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 194: Laos
$ws.Cells.Item(194, 1).Value = 'Laos'
$ws.Cells.Item(194, 3).Value = 0

# Row 195: Namibia
$ws.Cells.Item(195, 1).Value = 'Namibia'
$ws.Cells.Item(195, 3).Value = 1

# Row 209: Seychelles
$ws.Cells.Item(209, 1).Value = 'Seychelles'

# Row 210: Groenlandia
$ws.Cells.Item(210, 1).Value = 'Groenlandia'
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Montserrat
$ws.Cells.Item(211, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1

# Row 214: Sahara Occidental
$ws.Cells.Item(214, 1).Value = 'Sahara Occidental'

# Row 215: Bonaire, San Eustaquio y Saba
$ws.Cells.Item(215, 1).Value = 'Bonaire, San Eustaquio y Saba'
